# Rewrites the Questions sheet: replaces the old single-line, Python-repr-style
# "questions = [...]" shared string in A1 with a pretty-printed (json.dumps,
# indent=4) version, strips the stray trailing whitespace that was on two of
# the question titles, and removes row 2 (which used to duplicate-hold the
# string while A1 held a placeholder 0). Also drops A1's old bold+bordered
# style so it falls back to the workbook's default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = 'questions = [
    {
        "title": "You are a systems engineer working on an automation script. You collect the running services under a virtual machine VM1 using Get-Service | where {$_.Status -eq \"Running\"} cmdlet. You need to complete the above command to receive a JSON object as output.What should the command look like?",
        "ques_type": 2,
        "options": [
            "Get-Service | where {$_.Status -eq \"Running\"} | ConvertTo-Json",
            "Get-Service | where {$_.Status -eq \"Running\"} | ExportTo-Json",
            "Get-Service | where {$_.Status -eq \"Running\"} - json",
            "Get-Service | where {$_.Status -eq \"Running\"} -&gt json"
        ],
        "score": "Get-Service | where {$_.Status -eq \"Running\"} | ConvertTo-Json"
    },
    {
        "title": "You are an automation engineer working on a script to retrieve all the details about network cards in virtual machine VM1. You do not remember the exact command to do that, though you know it contains the word NetAdapter.How can you retrieve all available commands in PowerShell if you know only this part of the command?",
        "ques_type": 2,
        "options": [
            "Get-PSCommand *NetAdapter*",
            "Get-Command ?NetAdapter?",
            "Get-Command NetAdapter",
            "Get-Command *NetAdapter*"
        ],
        "score": "Get-Command *NetAdapter*"
    },
    {
        "title": "You are an Azure administrator, and you create a new variable named newvar in an Azure Automation account. You create a new PowerShell runbook, and you need to retrieve the value of the variable newvar inside runbook\u2019s code.Which command will retrieve the value?",
        "ques_type": 2,
        "options": [
            "$var = Get-AzAutomationVariable -Name ''newvar''",
            "$var = Get-AutomationVariable -Name ''newvar''",
            "$var = Get-Variable -Name ''newvar''",
            "$var = New-AzAutomationVariable -Name ''newvar''",
            "$var = New-Variable -Name \u2018newvar\u2019"
        ],
        "score": "$var = Get-AutomationVariable -Name ''newvar''"
    },
    {
        "title": "You are an Azure administrator. You receive a request to programmatically generate a report with the audit logs from Azure Active Directory. You write a PowerShell script, where you use the following as the API: https://graph.windows.net/{tenant}/activities/audit?api-version=beta Which PowerShell command should you use to call the API?",
        "ques_type": 2,
        "options": [
            "Invoke-WebRequest",
            "Get-API",
            "Invoke-API",
            "Get-WebRequest"
        ],
        "score": "Invoke-WebRequest"
    }
]'

# A1 previously carried the bold/centered/bordered style (s="1"); clear it so
# the cell reverts to the default style, matching the target file.
$ws.Range("A1").ClearFormats()

# Row 2 (which held the real shared-string text) is removed entirely; its
# content moves up into A1.
$ws.Range("A2").Delete()

$ws.Range("A1").Value = $questionsText

# Assigning a multi-line string makes the host auto-expand the row height;
# AutoFit() restores the implicit/default row height so no ht=/customHeight=
# attribute is written out, matching the target sheet XML.
$ws.Range("A1").EntireRow.AutoFit()
